$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.545.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.15%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.919.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.37%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.23%  "

# Row 6
$ws.Range("E6").Value = "  +0.00%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4800"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.42%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2899"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06720"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.45%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "111.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.65%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.15"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.78%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.914.56"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.50%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07565"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.62%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.309"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.85%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6727"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.55%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "299.11"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.04%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.528.78"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.25%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.594"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.70%  "

# Row 19
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.12%  "

# Row 20
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.02%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000007586"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.78%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.163.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.02%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9993"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.02%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.457"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.17%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.484"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.16%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.95%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.78%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.109"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.54%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1064"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.51%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.436"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.32%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.157"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.04%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.090"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.05%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05019"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.93%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7401"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.32%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.141"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.25%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.09%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.736"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.43%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02021"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.81%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.678"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.44%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "110.71"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.23%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.024"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.23%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4458"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.60%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8635"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.34%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.32%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.865"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.05%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.0000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.03%  "

# Row 47
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.260"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.34%  "

# Row 48
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "49.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.69%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.291"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.16%  "

# Row 50
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1239"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.04%  "

# Row 51
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2552"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.81%  "

